$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (current dbExcel/Neo4jData column),
# shifting the old B/C columns to C/D. Inserting this way carries the
# formatting (incl. the wrap-text style) of column A onto the new column B,
# same as row 2's A2 style, matching the target A2/B2 styling.
$ws.Columns("B").Insert()

# The freshly inserted column has no explicit width yet; give it the same
# display width as column A (75.8-ish characters).
$ws.Columns("B").ColumnWidth = 75

# New header cell for the stat-bar query column.
$ws.Range("B1").Value = "StatQuery"

# New stat-bar Neo4j query text (wrap-text formatting already carried over
# from the Insert above, matching A2's style).
$ws.Range("B2").Value = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Gordon Setter']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"
$ws.Range("B2").WrapText = $true

Write-Output "done"
